$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting (styles) of column R (rows 2-6) into the new column S
# so the new cells inherit the same number formats / borders / fonts as
# their row counterparts in column R.
$ws.Range("R2:R6").Copy()
$ws.Range("S2:S6").PasteSpecial(-4122)  # xlPasteFormats

# Populate the new column S values (row 2 stays blank, like R2)
$ws.Range("S3").Value = 2022
$ws.Range("S4").Value = 265803
$ws.Range("S5").Value = 3.8
$ws.Range("S6").Value = 33.6

# Update the active selection to match the recorded cursor position in the diff
$ws.Range("C19").Select()
